# Fix "Recorded By" (column G) attendance-audit entries.
#
# Several rows have a comma-separated list of recorder identities whose
# first entry is the literal "System" (e.g. "System, dnasr281@gmail.com").
# The canonical/expected ordering puts "System" last instead of first, so
# for every such cell we reverse the order of the comma-separated list
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
# Cells that are just "System" alone, or where "System" isn't the first
# token, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$modified = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.StartsWith("System,")) {
        $parts = $val -split ", "
        $n = $parts.Count

        $reversed = @()
        for ($i = $n - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i]
        }

        $cell.Value2 = $reversed -join ", "
        $modified++
    }
}

Write-Output "Reordered Recorded-By list in $modified row(s)."
